$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Global short-url rename (affects every data row via shared string)
$ws.Cells.Replace("9ffD9z", "p1u1FQ")

# 2. Country-of-origin text fields shift up (Palestinian entry removed, Q3 2025 refresh)
$ws.Range("G382").Value = "Libya"
$ws.Range("H382").Value = "LBY"
$ws.Range("I382").Value = "LBY"
$ws.Range("G383").Value = "Mali"
$ws.Range("H383").Value = "MLI"
$ws.Range("I383").Value = "MLI"
$ws.Range("G384").Value = "Niger"
$ws.Range("H384").Value = "NGR"
$ws.Range("I384").Value = "NER"
$ws.Range("G385").Value = "Nigeria"
$ws.Range("H385").Value = "NIG"
$ws.Range("I385").Value = "NGA"
$ws.Range("G386").Value = "Rwanda"
$ws.Range("H386").Value = "RWA"
$ws.Range("I386").Value = "RWA"
$ws.Range("G387").Value = "Somalia"
$ws.Range("H387").Value = "SOM"
$ws.Range("I387").Value = "SOM"
$ws.Range("G388").Value = "South Sudan"
$ws.Range("H388").Value = "SSD"
$ws.Range("I388").Value = "SSD"
$ws.Range("G389").Value = "Sudan"
$ws.Range("H389").Value = "SUD"
$ws.Range("I389").Value = "SDN"
$ws.Range("G390").Value = "Syrian Arab Rep."
$ws.Range("H390").Value = "SYR"
$ws.Range("I390").Value = "SYR"
$ws.Range("G391").Value = "Yemen"
$ws.Range("H391").Value = "YEM"
$ws.Range("I391").Value = "YEM"

# 3. Numeric-looking fields must stay stored as text, matching source data convention
$ws.Range("N373").NumberFormat = "@"
$ws.Range("N373").Value = "6"
$ws.Range("N375").NumberFormat = "@"
$ws.Range("N375").Value = "137683"
$ws.Range("O375").NumberFormat = "@"
$ws.Range("O375").Value = "4111"
$ws.Range("Q376").NumberFormat = "@"
$ws.Range("Q376").Value = "220610"
$ws.Range("T376").NumberFormat = "@"
$ws.Range("T376").Value = "317279"
$ws.Range("N377").NumberFormat = "@"
$ws.Range("N377").Value = "9059"
$ws.Range("O377").NumberFormat = "@"
$ws.Range("O377").Value = "68"
$ws.Range("P377").NumberFormat = "@"
$ws.Range("P377").Value = "17648"
$ws.Range("N378").NumberFormat = "@"
$ws.Range("N378").Value = "6"
$ws.Range("N379").NumberFormat = "@"
$ws.Range("N379").Value = "207"
$ws.Range("O379").NumberFormat = "@"
$ws.Range("O379").Value = "18"
$ws.Range("N380").NumberFormat = "@"
$ws.Range("N380").Value = "14"
$ws.Range("O380").NumberFormat = "@"
$ws.Range("O380").Value = "13"
$ws.Range("N381").NumberFormat = "@"
$ws.Range("N381").Value = "25"
$ws.Range("O381").NumberFormat = "@"
$ws.Range("O381").Value = "12"
$ws.Range("F382").NumberFormat = "@"
$ws.Range("F382").Value = "107"
$ws.Range("N382").NumberFormat = "@"
$ws.Range("N382").Value = "5"
$ws.Range("O382").NumberFormat = "@"
$ws.Range("O382").Value = "10"
$ws.Range("F383").NumberFormat = "@"
$ws.Range("F383").Value = "126"
$ws.Range("O383").NumberFormat = "@"
$ws.Range("O383").Value = "15"
$ws.Range("F384").NumberFormat = "@"
$ws.Range("F384").Value = "139"
$ws.Range("N384").NumberFormat = "@"
$ws.Range("N384").Value = "887"
$ws.Range("O384").NumberFormat = "@"
$ws.Range("O384").Value = "5"
$ws.Range("F385").NumberFormat = "@"
$ws.Range("F385").Value = "141"
$ws.Range("N385").NumberFormat = "@"
$ws.Range("N385").Value = "21484"
$ws.Range("O385").NumberFormat = "@"
$ws.Range("O385").Value = "398"
$ws.Range("P385").NumberFormat = "@"
$ws.Range("P385").Value = "214"
$ws.Range("F386").NumberFormat = "@"
$ws.Range("F386").Value = "161"
$ws.Range("N386").NumberFormat = "@"
$ws.Range("N386").Value = "5"
$ws.Range("O386").NumberFormat = "@"
$ws.Range("O386").Value = "0"
$ws.Range("F387").NumberFormat = "@"
$ws.Range("F387").Value = "172"
$ws.Range("N387").NumberFormat = "@"
$ws.Range("N387").Value = "5"
$ws.Range("O387").NumberFormat = "@"
$ws.Range("O387").Value = "0"
$ws.Range("F388").NumberFormat = "@"
$ws.Range("F388").Value = "179"
$ws.Range("N388").NumberFormat = "@"
$ws.Range("N388").Value = "75"
$ws.Range("O388").NumberFormat = "@"
$ws.Range("O388").Value = "0"
$ws.Range("P388").NumberFormat = "@"
$ws.Range("P388").Value = "0"
$ws.Range("F389").NumberFormat = "@"
$ws.Range("F389").Value = "177"
$ws.Range("N389").NumberFormat = "@"
$ws.Range("N389").Value = "1109357"
$ws.Range("O389").NumberFormat = "@"
$ws.Range("O389").Value = "3104"
$ws.Range("F390").NumberFormat = "@"
$ws.Range("F390").Value = "185"
$ws.Range("N390").NumberFormat = "@"
$ws.Range("N390").Value = "38"
$ws.Range("O390").NumberFormat = "@"
$ws.Range("O390").Value = "15"
$ws.Range("F391").NumberFormat = "@"
$ws.Range("F391").Value = "211"
$ws.Range("N391").NumberFormat = "@"
$ws.Range("N391").Value = "5"

# 4. Drop the trailing three rows (392:394) now superseded by the refreshed data
$ws.Rows("392:394").Delete()
